$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Row 2
$ws1.Range("D2").Value = 2
$ws1.Range("H2").Value = 6.4
$ws1.Range("L2").Value = 1.01

# Row 3
$ws1.Range("D3").Value = 2
$ws1.Range("H3").Value = 5.4
$ws1.Range("L3").Value = 0.9

# Row 4
$ws1.Range("H4").Value = 4.78
$ws1.Range("L4").Value = 1.18

# Row 5
$ws1.Range("H5").Value = 3.78
$ws1.Range("L5").Value = 1.03

# Row 6
$ws1.Range("H6").Value = 2.78
$ws1.Range("L6").Value = 1.14

# Row 7
$ws1.Range("H7").Value = 1.78
$ws1.Range("L7").Value = 1.02

# Row 8
$ws1.Range("H8").Value = 0.86
$ws1.Range("J8").Value = "Urgent"
$ws1.Range("L8").Value = 0.92

# Row 9
$ws1.Range("H9").Value = 0
$ws1.Range("I9").Value = "High"
$ws1.Range("L9").Value = 0.84

# Row 10
$ws1.Range("D10").Value = 2
$ws1.Range("L10").Value = 0.89

# Row 11
$ws1.Range("D11").Value = 2
$ws1.Range("L11").Value = 1.09

# Row 12
$ws1.Range("L12").Value = 1.01

# Row 13
$ws1.Range("L13").Value = 1.07

# Row 14
$ws1.Range("L14").Value = 1.17

# Row 15
$ws1.Range("L15").Value = 1.01

# Row 16
$ws1.Range("L16").Value = 0.89

# Row 17
$ws1.Range("L17").Value = 0.84

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

# Leading apostrophe forces these numeric-looking values to be stored
# as text, matching the original inlineStr cell type.
$ws2.Range("B9").Value = "'34"
$ws2.Range("B10").Value = "'18"
$ws2.Range("B11").Value = "'10"
